$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.867.41"
$ws.Range("E2").Value = "  +5.28%  "

$ws.Range("D3").Value = "1.878.52"
$ws.Range("E3").Value = "  +4.07%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "283.83"
$ws.Range("E5").Value = "  +2.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9986"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5213"
$ws.Range("E7").Value = "  +4.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3547"
$ws.Range("E8").Value = "  +1.46%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.44"
$ws.Range("E9").Value = "  +3.74%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07078"
$ws.Range("E10").Value = "  +6.23%  "

$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.35"
$ws.Range("E11").Value = "  +2.11%  "

$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8228"
$ws.Range("E12").Value = "  -1.32%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07752"
$ws.Range("E13").Value = "  -0.85%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.873.82"
$ws.Range("E14").Value = "  +3.79%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.205"
$ws.Range("E15").Value = "  +3.10%  "

$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.14"
$ws.Range("E16").Value = "  +3.46%  "

$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9987"
$ws.Range("E17").Value = "  -0.28%  "

$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.52"
$ws.Range("E18").Value = "  +4.53%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008172"
$ws.Range("E19").Value = "  +3.21%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9985"
$ws.Range("E20").Value = "  -0.33%  "

$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "26.902.72"
$ws.Range("E21").Value = "  +5.06%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.792"
$ws.Range("E22").Value = "  +1.73%  "

$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.22"
$ws.Range("E23").Value = "  +2.76%  "

$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.254"
$ws.Range("E24").Value = "  +3.33%  "

$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.434"
$ws.Range("E25").Value = "  +15.20%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.82"
$ws.Range("E26").Value = "  +3.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.50"
$ws.Range("E27").Value = "  +3.29%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.667"
$ws.Range("E28").Value = "  +0.75%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.83"
$ws.Range("E29").Value = "  +2.87%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.425"
$ws.Range("E30").Value = "  +3.08%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.384"
$ws.Range("E31").Value = "  +4.27%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08885"
$ws.Range("E32").Value = "  +0.27%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04942"
$ws.Range("E33").Value = "  +2.78%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.187"
$ws.Range("E34").Value = "  +5.75%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7511"
$ws.Range("E35").Value = "  +2.81%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.873"
$ws.Range("E36").Value = "  +1.10%  "

$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.298"
$ws.Range("E37").Value = "  +9.03%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.422"
$ws.Range("E38").Value = "  +7.20%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5361"
$ws.Range("E39").Value = "  +3.52%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01888"
$ws.Range("E40").Value = "  +1.76%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9805"
$ws.Range("E41").Value = "  +2.23%  "

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "116.93"
$ws.Range("E42").Value = "  +3.66%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.321"
$ws.Range("E43").Value = "  +2.53%  "

$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.266"
$ws.Range("E44").Value = "  +2.85%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4659"
$ws.Range("E45").Value = "  +1.80%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9979"
$ws.Range("E46").Value = "  -0.22%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1374"
$ws.Range("E47").Value = "  -0.35%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.477"
$ws.Range("E48").Value = "  +3.02%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.83"
$ws.Range("E49").Value = "  +3.29%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.524"
$ws.Range("E50").Value = "  +2.51%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05940"
$ws.Range("E51").Value = "  +2.24%  "
